$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param([string]$CellRef, [string]$Value)
    $c = $ws.Range($CellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $Value
    $c.Style = $origStyle
}

Set-CellText "D2" "20.546.85"
Set-CellText "E2" "  +1.58%  "
Set-CellText "D3" "1.473.22"
Set-CellText "E3" "  +2.04%  "
Set-CellText "E4" "  +0.05%  "
Set-CellText "D5" "0.9570"
Set-CellText "E5" "  +4.94%  "
Set-CellText "D6" "277.33"
Set-CellText "E6" "  -0.27%  "
Set-CellText "D7" "0.3611"
Set-CellText "E7" "  -1.29%  "
Set-CellText "D8" "0.3077"
Set-CellText "E8" "  -1.44%  "
Set-CellText "D9" "39.66"
Set-CellText "E9" "  +1.25%  "
Set-CellText "E10" "  +4.59%  "
Set-CellText "D11" "0.06661"
Set-CellText "E11" "  +1.82%  "
Set-CellText "D12" "1.001"
Set-CellText "E12" "  +0.04%  "
Set-CellText "D13" "5.528"
Set-CellText "E13" "  +2.49%  "
Set-CellText "D14" "18.15"
Set-CellText "E14" "  +2.86%  "
Set-CellText "D15" "6.178"
Set-CellText "E15" "  +1.89%  "
Set-CellText "D16" "0.9566"
Set-CellText "E16" "  +1.52%  "
Set-CellText "D17" "0.00001028"
Set-CellText "E17" "  +1.16%  "
Set-CellText "D18" "1.472.03"
Set-CellText "E18" "  +1.85%  "
Set-CellText "D19" "0.05930"
Set-CellText "E19" "  +5.20%  "
Set-CellText "D20" "68.93"
Set-CellText "E20" "  +1.06%  "
Set-CellText "D21" "5.493"
Set-CellText "E21" "  +1.71%  "
Set-CellText "D22" "14.53"
Set-CellText "E22" "  +0.68%  "
Set-CellText "E23" "  +2.78%  "
Set-CellText "D24" "2.262"
Set-CellText "E24" "  +0.62%  "
Set-CellText "D25" "20.547.05"
Set-CellText "E25" "  +1.57%  "
Set-CellText "D26" "143.27"
Set-CellText "E26" "  +3.90%  "
Set-CellText "D27" "2.127"
Set-CellText "E27" "  -2.20%  "
Set-CellText "E28" "  +0.97%  "
Set-CellText "D29" "1.631.62"
Set-CellText "E29" "  +2.16%  "
Set-CellText "D30" "113.90"
Set-CellText "E30" "  +3.15%  "
Set-CellText "D31" "3.904"
Set-CellText "E31" "  +2.22%  "
Set-CellText "D32" "4.974"
Set-CellText "E32" "  +2.71%  "
Set-CellText "D33" "0.08016"
Set-CellText "E33" "  +4.16%  "
Set-CellText "D34" "0.8061"
Set-CellText "E34" "  +0.31%  "
Set-CellText "D35" "1.514"
Set-CellText "E35" "  +3.67%  "
Set-CellText "D36" "1.217"
Set-CellText "E36" "  +6.34%  "
Set-CellText "D37" "0.05746"
Set-CellText "E37" "  -3.48%  "
Set-CellText "D38" "4.730"
Set-CellText "E38" "  +0.86%  "
Set-CellText "D39" "0.02058"
Set-CellText "E39" "  +3.21%  "
Set-CellText "B40" "Aptos"
Set-CellText "C40" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-CellText "D40" "10.40"
Set-CellText "E40" "  +2.24%  "
Set-CellText "B41" "Frax"
Set-CellText "C41" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-CellText "D41" "0.9578"
Set-CellText "E41" "  +3.01%  "
Set-CellText "D42" "0.1872"
Set-CellText "E42" "  +1.65%  "
Set-CellText "E43" "  +4.83%  "
Set-CellText "D44" "0.5276"
Set-CellText "E44" "  +0.72%  "
Set-CellText "D45" "3.520"
Set-CellText "E45" "  -0.13%  "
Set-CellText "D46" "12.16"
Set-CellText "E46" "  +0.36%  "
Set-CellText "D47" "118.41"
Set-CellText "E47" "  -0.43%  "
Set-CellText "D48" "0.5207"
Set-CellText "E48" "  +1.19%  "
Set-CellText "D49" "1.816"
Set-CellText "E49" "  +3.26%  "
Set-CellText "D50" "0.06481"
Set-CellText "E50" "  +2.22%  "
Set-CellText "D51" "0.9868"
Set-CellText "E51" "  -0.51%  "
